$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.057.06"
$ws.Range("E2").Value = "  -0.39%  "

$ws.Range("D3").Value = "2.589.10"
$ws.Range("E3").Value = "  +0.74%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.20%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("E8").Value = "  -2.66%  "

$ws.Range("E9").Value = "  -1.69%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.60%  "

$ws.Range("E11").Value = "  -0.62%  "

$ws.Range("E12").Value = "  -1.71%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.54%  "

$ws.Range("D14").Value = "3.048.86"
$ws.Range("E14").Value = "  +0.59%  "

$ws.Range("D15").Value = "62.972.09"
$ws.Range("E15").Value = "  -0.48%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000146"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.87%  "

$ws.Range("D17").Value = "2.581.35"
$ws.Range("E17").Value = "  +0.23%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.93%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "341.90"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.39%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.48%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.12%  "

$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.42%  "

$ws.Range("E25").Value = "  -1.84%  "

$ws.Range("E26").Value = "  -3.19%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.27%  "

$ws.Range("E28").Value = "  -0.29%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.05%  "

$ws.Range("E30").Value = "  -2.82%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "464.58"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.81%  "

$ws.Range("D32").Value = "0.0₃0803"
$ws.Range("E32").Value = "  -2.87%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.68"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.63%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "176.80"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.00%  "

$ws.Range("E35").Value = "  +0.07%  "

$ws.Range("E36").Value = "  -2.42%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.93"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.81%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.51"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.05%  "

$ws.Range("E40").Value = "  -3.11%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "159.38"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.00%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "40.10"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.97%  "

$ws.Range("E43").Value = "  -2.99%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.68%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.638"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.60%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0540"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.74%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0963"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.25%  "

$ws.Range("E48").Value = "  -1.74%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.03"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.51%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.40"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.69"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.17%  "
